$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source sheet stores every data cell as text (inline strings / shared
# strings), including values that look numeric (prices, hour, etc.). To
# reproduce that faithfully through Excel COM we temporarily force a Text
# number format before assigning the value (otherwise Excel would silently
# coerce strings like "25.10" into the number 25.1) and then restore the
# default "Normal" style afterwards so no stray per-cell style/quote-prefix
# marker is left behind, matching the original (unstyled) cells exactly.
$updates = @(
    @{ Cell = 'D2'; Value = '244.93' },
    @{ Cell = 'G2'; Value = '20' },
    @{ Cell = 'D3'; Value = '25.10' },
    @{ Cell = 'G3'; Value = '20' },
    @{ Cell = 'D4'; Value = '4.992' },
    @{ Cell = 'G4'; Value = '20' },
    @{ Cell = 'D5'; Value = '0.05609' },
    @{ Cell = 'G5'; Value = '20' },
    @{ Cell = 'D6'; Value = '6.576' },
    @{ Cell = 'G6'; Value = '20' },
    @{ Cell = 'D7'; Value = '3.009' },
    @{ Cell = 'G7'; Value = '20' },
    @{ Cell = 'D8'; Value = '0.8093' },
    @{ Cell = 'G8'; Value = '20' },
    @{ Cell = 'D9'; Value = '0.8413' },
    @{ Cell = 'G9'; Value = '20' },
    @{ Cell = 'D10'; Value = '0.1338' },
    @{ Cell = 'G10'; Value = '20' },
    @{ Cell = 'B11'; Value = 'MandalaExchangeToken' },
    @{ Cell = 'C11'; Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx' },
    @{ Cell = 'D11'; Value = '0.06942' },
    @{ Cell = 'E11'; Value = '10MandalaExchangeTokenMDX' },
    @{ Cell = 'G11'; Value = '20' },
    @{ Cell = 'D12'; Value = '0.02838' },
    @{ Cell = 'G12'; Value = '20' },
    @{ Cell = 'D13'; Value = '0.09405' },
    @{ Cell = 'G13'; Value = '20' },
    @{ Cell = 'D14'; Value = '0.001513' },
    @{ Cell = 'G14'; Value = '20' },
    @{ Cell = 'D15'; Value = '0.0005953' },
    @{ Cell = 'G15'; Value = '20' },
    @{ Cell = 'D16'; Value = '0.006135' },
    @{ Cell = 'G16'; Value = '20' },
    @{ Cell = 'D17'; Value = '3.497' },
    @{ Cell = 'G17'; Value = '20' },
    @{ Cell = 'D18'; Value = '2.092' },
    @{ Cell = 'G18'; Value = '20' },
    @{ Cell = 'G19'; Value = '20' },
    @{ Cell = 'B20'; Value = 'LiechtensteinCryptoassetsExchange' },
    @{ Cell = 'C20'; Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx' },
    @{ Cell = 'D20'; Value = '0.03266' },
    @{ Cell = 'E20'; Value = '19LiechtensteinCryptoassetsExchangeLCX' },
    @{ Cell = 'G20'; Value = '20' },
    @{ Cell = 'D21'; Value = '0.1292' },
    @{ Cell = 'G21'; Value = '20' },
    @{ Cell = 'D22'; Value = '3.749' },
    @{ Cell = 'G22'; Value = '20' },
    @{ Cell = 'D23'; Value = '0.04672' },
    @{ Cell = 'G23'; Value = '20' },
    @{ Cell = 'G24'; Value = '20' },
    @{ Cell = 'D25'; Value = '0.001243' },
    @{ Cell = 'G25'; Value = '20' },
    @{ Cell = 'D26'; Value = '0.004525' },
    @{ Cell = 'G26'; Value = '20' },
    @{ Cell = 'D27'; Value = '0.00009696' },
    @{ Cell = 'E27'; Value = '26NitroExNTXBestin24h' },
    @{ Cell = 'G27'; Value = '20' },
    @{ Cell = 'D28'; Value = '0.0001939' },
    @{ Cell = 'G28'; Value = '20' },
    @{ Cell = 'G29'; Value = '20' },
    @{ Cell = 'G30'; Value = '20' },
    @{ Cell = 'G31'; Value = '20' },
    @{ Cell = 'G32'; Value = '20' },
    @{ Cell = 'G33'; Value = '20' },
    @{ Cell = 'G34'; Value = '20' },
    @{ Cell = 'G35'; Value = '20' },
    @{ Cell = 'G36'; Value = '20' },
    @{ Cell = 'G37'; Value = '20' },
    @{ Cell = 'G38'; Value = '20' },
    @{ Cell = 'G39'; Value = '20' },
    @{ Cell = 'D40'; Value = '0.03646' },
    @{ Cell = 'G40'; Value = '20' },
    @{ Cell = 'D41'; Value = '0.1348' },
    @{ Cell = 'G41'; Value = '20' },
    @{ Cell = 'D42'; Value = '0.006237' },
    @{ Cell = 'E42'; Value = '41KickTokenKICK' },
    @{ Cell = 'G42'; Value = '20' },
    @{ Cell = 'D43'; Value = '0.002530' },
    @{ Cell = 'G43'; Value = '20' },
    @{ Cell = 'D44'; Value = '0.008069' },
    @{ Cell = 'G44'; Value = '20' },
    @{ Cell = 'D45'; Value = '0.00005278' },
    @{ Cell = 'G45'; Value = '20' },
    @{ Cell = 'G46'; Value = '20' },
    @{ Cell = 'D47'; Value = '0.1799' },
    @{ Cell = 'G47'; Value = '20' },
    @{ Cell = 'D48'; Value = '0.002042' },
    @{ Cell = 'G48'; Value = '20' },
    @{ Cell = 'D49'; Value = '0.00002099' },
    @{ Cell = 'G49'; Value = '20' },
    @{ Cell = 'D50'; Value = '0.0001999' },
    @{ Cell = 'G50'; Value = '20' },
    @{ Cell = 'G51'; Value = '20' }
)

foreach ($update in $updates) {
    $cell = $ws.Range($update.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $update.Value
    $cell.Style = "Normal"
}
